# Create Initial route files
#
# Removes the stray empty, numbered ("ListParagraph") bullet that was left
# right after the "Create Map => '/map-form" bullet (just before the blank
# paragraph that precedes "POST Requests"). The bullet has no text and
# carries only list formatting (pStyle=ListParagraph + numPr), so find it
# by that combination rather than a hard-coded paragraph index, and delete
# the whole paragraph (its Range, including the paragraph mark).

$d = $word.ActiveDocument

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text.Trim()
    if ($text -eq "" -and $p.Range.ListFormat.ListType -ne 0) {
        $p.Range.Delete()
    }
}
